$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 338 corresponds to "na?o informado" which needs to be removed entirely.
# Deleting the entire row shifts all subsequent rows up by one.
$ws.Rows.Item(338).Delete()
